# InterimPricing.xlsx — "added feed till september"
# Append two new monthly observations (Aug-2022 and Sep-2022) to the
# Date/Value feed on Sheet1, right after the existing last row (104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateFormat = "[$-409]mmm\-yy;@"

# Row 105: 2022-08-01 -> 686
$ws.Cells.Item(105, 1).Value = 44774
$ws.Cells.Item(105, 1).NumberFormat = $dateFormat
$ws.Cells.Item(105, 2).Value = 686

# Row 106: 2022-09-01 -> 675
$ws.Cells.Item(106, 1).Value = 44805
$ws.Cells.Item(106, 1).NumberFormat = $dateFormat
$ws.Cells.Item(106, 2).Value = 675

# Reflect the final view/selection state after the edit.
$ws.Range("B107").Select()
